$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.50%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.06%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.121"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.96%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07370"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.45%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.705"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.09%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.640"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "13.60%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9189"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.53%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1673"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.11%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07112"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.85%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07943"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.05%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02977"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.09%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09907"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.81%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006145"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.61%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.66%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.18%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.87%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1332"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.14%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.553"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.41%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04617"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.79%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.39%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.97%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004423"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.29%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.58%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.80%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01681"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04412"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.46%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007174"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.13%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1327"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.00%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002135"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.08%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01103"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005991"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.60%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.928"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.87%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-36.88%"
